$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values per row (columns B, C, D, F, G, H, I, L)
# row number -> @{ Col = Value; ... }
$updates = @{
    2 = @{ "B"=17.2042400402971; "C"=10.76154531977254; "D"=4.794905354512239; "F"=25.92069599919177; "G"=31.9694579864587; "H"=14.76500625785558; "I"=22.48305991620906; "L"=10.68012773888499 }
    3 = @{ "B"=16.5505715605408; "C"=10.38812948220805; "D"=4.796473984247393; "F"=25.8765808220838; "G"=31.83880975046162; "H"=14.81381235575958; "I"=22.61094462936199; "L"=10.65435446621486 }
    4 = @{ "B"=16.13791824122116; "C"=10.15018367922485; "D"=4.797787154373998; "F"=25.86002218862061; "G"=31.77475749132969; "H"=14.84783800416863; "I"=22.6964944424178; "L"=10.64077390943622 }
    5 = @{ "B"=15.96716950342466; "C"=10.05113272247192; "D"=4.798410698506605; "F"=25.85591990297197; "G"=31.75273158593248; "H"=14.86271966201906; "I"=22.73311473546034; "L"=10.63580763775287 }
    6 = @{ "B"=15.93866831745522; "C"=10.03456231899975; "D"=4.798519589402141; "F"=25.85539837689981; "G"=31.74932052584073; "H"=14.86525197136448; "I"=22.73930142018326; "L"=10.63501738637491 }
    7 = @{ "B"=16.13562559480692; "C"=10.14885615912752; "D"=4.79779520520331; "F"=25.85995615754067; "G"=31.77444393119771; "H"=14.84803459622916; "I"=22.69698121046391; "L"=10.6407046287864 }
    8 = @{ "B"=16.98134292768366; "C"=10.63464317073607; "D"=4.795373808911795; "F"=25.90329986272962; "G"=31.92106419440368; "H"=14.78098953345207; "I"=22.52569033781681; "L"=10.67077759356301 }
    9 = @{ "B"=18.54041936997608; "C"=11.51491300198477; "D"=4.793385703741165; "F"=26.07175208380803; "G"=32.33597490846429; "H"=14.68191586532678; "I"=22.24596712341786; "L"=10.74737941081023 }
    10 = @{ "B"=19.61403132591442; "C"=12.11339225831947; "D"=4.793585630430098; "F"=26.24603879760985; "G"=32.71670806317371; "H"=14.62914032585731; "I"=22.07526812655732; "L"=10.81413371437584 }
    11 = @{ "B"=20.0849282221725; "C"=12.37447253764065; "D"=4.794032519282863; "F"=26.33614590986724; "G"=32.90585224413311; "H"=14.60952882461049; "I"=22.00529411051222; "L"=10.84670801024137 }
    12 = @{ "B"=20.26059653847699; "C"=12.47168147938005; "D"=4.794252529651604; "F"=26.37180516257298; "G"=32.97971300238474; "H"=14.60273851453258; "I"=21.97991096974274; "L"=10.85935388845937 }
    13 = @{ "B"=20.22288303516772; "C"=12.45082019458812; "D"=4.794202894249958; "F"=26.3640572520426; "G"=32.96370730190132; "H"=14.60417257665409; "I"=21.98532794944512; "L"=10.85661666587056 }
    14 = @{ "B"=20.09943437769376; "C"=12.38250343767251; "D"=4.794049604189192; "F"=26.33904892085985; "G"=32.91188435212143; "H"=14.60895741081926; "I"=22.00318342824478; "L"=10.84774220753759 }
    15 = @{ "B"=20.02346971006827; "C"=12.34044035424786; "D"=4.793962311200232; "F"=26.32393022801732; "G"=32.88043067726446; "H"=14.61197121691414; "I"=22.01426587506848; "L"=10.84234660029165 }
    16 = @{ "B"=19.58289315671128; "C"=12.09610081975826; "D"=4.793563553372341; "F"=26.24036635584973; "G"=32.70466320304239; "H"=14.63051079289942; "I"=22.07999648299284; "L"=10.812048762287 }
    17 = @{ "B"=19.3080313420214; "C"=11.94330743494571; "D"=4.793409808123632; "F"=26.1918625017568; "G"=32.6008813872437; "H"=14.6430131245315; "I"=22.122293984976; "L"=10.79402275836244 }
    18 = @{ "B"=19.14829821388992; "C"=11.8543757993085; "D"=4.793354904718239; "F"=26.1649844462503; "G"=32.54269368747559; "H"=14.65061772050298; "I"=22.1473444375955; "L"=10.7838628842278 }
    19 = @{ "B"=19.09393807124466; "C"=11.82408656176488; "D"=4.793342085837509; "F"=26.15605971105693; "G"=32.5232524231465; "H"=14.65326341520763; "I"=22.15594978454957; "L"=10.78045888237112 }
    20 = @{ "B"=19.33746160529537; "C"=11.95968149554814; "D"=4.793422707609629; "F"=26.19692037037465; "G"=32.61177378490773; "H"=14.64163939617552; "I"=22.11771653899738; "L"=10.79592015593305 }
    21 = @{ "B"=20.1357671450899; "C"=12.40261504821184; "D"=4.794093254140639; "F"=26.34635290900433; "G"=32.92704581369543; "H"=14.60753469418422; "I"=21.99790851093186; "L"=10.85034047466494 }
    22 = @{ "B"=20.6420091938362; "C"=12.68242376311279; "D"=4.794827313409378; "F"=26.45296719320691; "G"=33.14609615276971; "H"=14.58895485780667; "I"=21.92610803207482; "L"=10.88771485329779 }
    23 = @{ "B"=20.37327601527442; "C"=12.53398461099991; "D"=4.794408596970393; "F"=26.39525315171446; "G"=33.02801576014595; "H"=14.59853060366732; "I"=21.9638310124655; "L"=10.86760439235915 }
    24 = @{ "B"=19.32416150353886; "C"=11.95228216993375; "D"=4.793416771400801; "F"=26.19463056847851; "G"=32.60684472435317; "H"=14.6422591606418; "I"=22.11978371970365; "L"=10.79506170807317 }
    25 = @{ "B"=18.13054162142456; "C"=11.28498500046825; "D"=4.793630324220357; "F"=26.01726866133058; "G"=32.21024272093796; "H"=14.70522134060632; "I"=22.31556708823188; "L"=10.72479532689341 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
